$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "SEPTIEMBRE    2 0 2 1   " (3rd sheet) - add an adjustments block
# (J33:L37), clear the now-redundant P33:P36 formulas down to their cached
# 0 value, and fill in the two manual closing figures F54 / F57.
# ---------------------------------------------------------------------------
$wsSept = $wb.Worksheets.Item(3)

$wsSept.Range("J33").Value = "SEPT-,21"
$wsSept.Range("K33").Value = "BATAS"
$wsSept.Range("L33").Value = 3422

$wsSept.Range("J34").Value = "SEPT-,21"
$wsSept.Range("K34").Value = "XXXXX"
$wsSept.Range("L34").Value = 4999.6000000000004

$wsSept.Range("J35").Value = "SEPT-,21"
$wsSept.Range("K35").Value = "XXXXX"
$wsSept.Range("L35").Value = 1195.68

$wsSept.Range("J36").Value = "SEPT-,21"
$wsSept.Range("K36").Value = "FUMIGACION"
$wsSept.Range("L36").Value = 1392

$wsSept.Range("J37").Value = "SEPT-,21"
$wsSept.Range("K37").Value = "ADT"
$wsSept.Range("L37").Value = 836.84

# These used to mirror N+M+L+I+C via a shared formula; they are now pinned
# to their previous cached value (0) so the new L-column entries above do
# not flow into the P/Q cuadre columns.
$wsSept.Range("P33").Value = 0
$wsSept.Range("P34").Value = 0
$wsSept.Range("P35").Value = 0
$wsSept.Range("P36").Value = 0

# Manual closing figures.
$wsSept.Range("F54").Value = -1424333.95
$wsSept.Range("F57").Value = 5704

[void]$wsSept.Activate()
$wsSept.Range("K34").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "REMISIONES  SEPTIEMBRE  2021  " (4th sheet) - log the remaining
# "B" remittances for the period (rows 28-44).
# ---------------------------------------------------------------------------
$wsRem = $wb.Worksheets.Item(4)

$wsRem.Range("B28").Value = "20629 B"

$wsRem.Range("A29").Value = 44460
$wsRem.Range("B29").Value = "20728 B"
$wsRem.Range("C29").Value = 36425.199999999997

$wsRem.Range("A30").Value = 44461
$wsRem.Range("B30").Value = "20764 B"
$wsRem.Range("C30").Value = 50568.800000000003

$wsRem.Range("A31").Value = 44462
$wsRem.Range("B31").Value = "20875 B"
$wsRem.Range("C31").Value = 59040.6

$wsRem.Range("A32").Value = 44463
$wsRem.Range("B32").Value = "21070 B"
$wsRem.Range("C32").Value = 79386.06

$wsRem.Range("A33").Value = 44463
$wsRem.Range("B33").Value = "21071 B"
$wsRem.Range("C33").Value = 975

$wsRem.Range("A34").Value = 44464
$wsRem.Range("B34").Value = "21130 B"
$wsRem.Range("C34").Value = 26659.84

$wsRem.Range("A35").Value = 44464
$wsRem.Range("B35").Value = "21206 B"
$wsRem.Range("C35").Value = 63162.2

$wsRem.Range("A36").Value = 44466
$wsRem.Range("B36").Value = "21336 B"
$wsRem.Range("C36").Value = 42269.1

$wsRem.Range("A37").Value = 44466
$wsRem.Range("B37").Value = "21390 B"
$wsRem.Range("C37").Value = 3041

$wsRem.Range("A38").Value = 44467
$wsRem.Range("B38").Value = "21463 B"
$wsRem.Range("C38").Value = 59570.38
$wsRem.Range("D38").Value = 44470
$wsRem.Range("E38").Value = 400000

$wsRem.Range("A39").Value = 44467
$wsRem.Range("B39").Value = "21474 B"
$wsRem.Range("C39").Value = 7939.6

$wsRem.Range("A40").Value = 44468
$wsRem.Range("B40").Value = "21559 B"
$wsRem.Range("C40").Value = 38874.400000000001

$wsRem.Range("A41").Value = 44469
$wsRem.Range("B41").Value = "21666 B"
$wsRem.Range("C41").Value = 92182.8

$wsRem.Range("A42").Value = 44470
$wsRem.Range("B42").Value = "21845 B"
$wsRem.Range("C42").Value = 66246

$wsRem.Range("A43").Value = 44471
$wsRem.Range("B43").Value = "21945 B"
$wsRem.Range("C43").Value = 85535.7

$wsRem.Range("A44").Value = 44471
$wsRem.Range("B44").Value = "21974 B"
$wsRem.Range("C44").Value = 1725

[void]$wsRem.Activate()
$wsRem.Range("A45").Select() | Out-Null

# ---------------------------------------------------------------------------
# Restore the sheet that should be active when the workbook re-opens.
# ---------------------------------------------------------------------------
[void]$wsSept.Activate()
